$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New header column I -> "Vasco" (must be last new shared string so it gets
# appended after the row labels below, matching author's entry order)

# Row 3 : Japon_death
$ws.Range("A3").Value = "Japon_death"
$ws.Range("B3").Value = 99
$ws.Range("C3").Value = 85
$ws.Range("D3").Value = 115
$ws.Range("F3").Value = 200
$ws.Range("I3").Value = 95

# Row 2 : Japon_total
$ws.Range("A2").Value = "Japon_total"
$ws.Range("B2").Value = 3900
$ws.Range("C2").Value = 4000
$ws.Range("D2").Value = 3800
$ws.Range("F2").Value = 6000
$ws.Range("I2").Value = 4800

# Row 4 : Chile_total
$ws.Range("A4").Value = "Chile_total"
$ws.Range("B4").Value = 7200
$ws.Range("C4").Value = 5000
$ws.Range("D4").Value = 7450
$ws.Range("F4").Value = 7000
$ws.Range("I4").Value = 6400

# Row 5 : Chile_death
$ws.Range("A5").Value = "Chile_death"
$ws.Range("B5").Value = 108
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 84
$ws.Range("F5").Value = 150
$ws.Range("I5").Value = 55

# Row 6 : Israel_total
$ws.Range("A6").Value = "Israel_total"
$ws.Range("B6").Value = 11100
$ws.Range("C6").Value = 11000
$ws.Range("D6").Value = 11250
$ws.Range("F6").Value = 12000
$ws.Range("I6").Value = 11000

# Row 7 : Israel_death
$ws.Range("A7").Value = "Israel_death"
$ws.Range("B7").Value = 88
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = 87
$ws.Range("F7").Value = 400
$ws.Range("H7").Value = $ws.Range("H6").Value2
$ws.Range("H6").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I7").Value = 75

# New header cell for the "Vasco" column
$ws.Range("I1").Value = "Vasco"

# Update the active selection to match the author's final cursor position
[void]$ws.Range("F7").Select()
